$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order (2005 .. 2306) for column E, rows 16-53.
$periods = @(
    "2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
    "2301","2302","2303","2304","2305","2306"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Valor Mora (column F) is 35112 for every period except the most recent one (2306),
# which keeps the 25749 value. The special value moved from row 16 to row 53 because
# the period list direction reversed.
$ws.Range("F16").Value = 35112
$ws.Range("F53").Value = 25749
